$wb = $excel.ActiveWorkbook

# New rows of game data (rows 5-7) appended identically to both sheets
$newRows = @(
    @("Glib", "Atari 2600", "0.0", "Unknown", "Selchow & Righter", "https://www.vgchartz.com/games/boxart/1137241ccc.jpg", "01/01/1983"),
    @("Phaser Patrol", "Atari 2600", "0.0", "Unknown", "Arcadia", "https://www.vgchartz.com/games/boxart/5877491ccc.jpg", "01/01/1982"),
    @("Kool-Aid Man", "Atari 2600", "0.0", "Mattel Interactive", "Mattel", "https://www.vgchartz.com/games/boxart/775863ccc.jpg", "01/01/1983")
)

foreach ($sheetName in @("games", "games_2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 5
    foreach ($row in $newRows) {
        $ws.Range("A" + $r).Value = $row[0]
        $ws.Range("B" + $r).Value = $row[1]

        # Units(m) column: keep as literal text like "0.0" instead of being
        # coerced to the number 0.
        $ws.Range("C" + $r).NumberFormat = "@"
        $ws.Range("C" + $r).Value = $row[2]
        $ws.Range("C" + $r).NumberFormat = "General"

        $ws.Range("D" + $r).Value = $row[3]
        $ws.Range("E" + $r).Value = $row[4]
        $ws.Range("F" + $r).Value = $row[5]

        # Release Date column: stored as plain text (matches existing rows'
        # text-formatted dates), not an auto-converted date serial.
        $ws.Range("G" + $r).NumberFormat = "@"
        $ws.Range("G" + $r).Value = $row[6]

        $r = $r + 1
    }
}

# Switch the active/selected sheet from "games_2" to "games", and update
# each sheet's remembered selection.
$wsGames = $wb.Worksheets.Item("games")
$wsGames2 = $wb.Worksheets.Item("games_2")

$wsGames2.Activate()
$wsGames2.Range("A6").Select()

$wsGames.Activate()
$wsGames.Range("A20").Select()
